$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 116.8660278320312
$ws.Cells.Item(3, 2).Value = 119.8711471557617
$ws.Cells.Item(4, 2).Value = 118.3105239868164
$ws.Cells.Item(5, 2).Value = 112.2388229370117
$ws.Cells.Item(6, 2).Value = 114.8602294921875
$ws.Cells.Item(7, 2).Value = 111.3124694824219
$ws.Cells.Item(8, 2).Value = 113.2200775146484
$ws.Cells.Item(9, 2).Value = 109.5124816894531
$ws.Cells.Item(10, 2).Value = 112.6587677001953
$ws.Cells.Item(11, 2).Value = 113.4131698608398
$ws.Cells.Item(12, 2).Value = 113.5954818725586
$ws.Cells.Item(13, 2).Value = 121.4106369018555
$ws.Cells.Item(14, 2).Value = 126.5259399414062
$ws.Cells.Item(15, 2).Value = 131.0348815917969
$ws.Cells.Item(16, 2).Value = 143.4014282226562
$ws.Cells.Item(17, 2).Value = 169.0320739746094
$ws.Cells.Item(18, 2).Value = 158.4875183105469
$ws.Cells.Item(19, 2).Value = 168.3096313476562
$ws.Cells.Item(20, 2).Value = 167.9345855712891
$ws.Cells.Item(21, 2).Value = 170.9794006347656
$ws.Cells.Item(22, 2).Value = 171.5580596923828
$ws.Cells.Item(23, 2).Value = 171.615234375
$ws.Cells.Item(24, 2).Value = 170.3724670410156
$ws.Cells.Item(25, 2).Value = 169.9477233886719
$ws.Cells.Item(26, 2).Value = 165.2580108642578
$ws.Cells.Item(27, 2).Value = 168.7827301025391
$ws.Cells.Item(28, 2).Value = 168.6980133056641
$ws.Cells.Item(29, 2).Value = 165.6544189453125
$ws.Cells.Item(30, 2).Value = 167.6476135253906
$ws.Cells.Item(31, 2).Value = 167.7669677734375
$ws.Cells.Item(32, 2).Value = 176.078369140625
$ws.Cells.Item(33, 2).Value = 190.8202209472656
$ws.Cells.Item(34, 2).Value = 178.8159332275391
$ws.Cells.Item(35, 2).Value = 218.2801513671875
$ws.Cells.Item(36, 2).Value = 222.4234161376953
$ws.Cells.Item(37, 2).Value = 222.9360809326172
$ws.Cells.Item(38, 2).Value = 205.1527252197266
$ws.Cells.Item(39, 2).Value = 194.72900390625
$ws.Cells.Item(40, 2).Value = 184.1470642089844
$ws.Cells.Item(41, 2).Value = 172.6604919433594
$ws.Cells.Item(42, 2).Value = 164.2326354980469
$ws.Cells.Item(43, 2).Value = 149.0408630371094
$ws.Cells.Item(44, 2).Value = 156.6593627929688
$ws.Cells.Item(45, 2).Value = 143.1354370117188
$ws.Cells.Item(46, 2).Value = 146.0276489257812
$ws.Cells.Item(47, 2).Value = 138.4981842041016
$ws.Cells.Item(48, 2).Value = 137.1892242431641
$ws.Cells.Item(49, 2).Value = 137.0011596679688
